$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos: value changes to the "Herlandí" bio line that used
#     to live at row 13 (old layout had the docente name one row below its
#     label; new layout moves it up under "Objetivos:" AND keeps a second
#     copy under "Método:" at row 18) ---
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

# --- Row 13: gains an A13 label ("Programa resumido:") and its B/C value
#     becomes "Semestral" ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# --- Row 14: label becomes "Short syllabus:", value unchanged text ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "1. Management Practice Areas. 2. Organizational structure"
$ws.Range("C14").Value = "1. Management Practice Areas. 2. Organizational structure"
$ws.Rows(14).RowHeight = 60

# --- Row 15: label becomes "Programa:", value becomes "01/01/2021" ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"
$ws.Rows(15).RowHeight = 120

# --- Row 16: label becomes "Syllabus:", value unchanged text ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings."
$ws.Range("C16").Value = "1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings."
$ws.Rows(16).RowHeight = 120

# --- Row 17: label becomes "Avaliação:" and loses its B/C value entirely ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Clear()
$ws.Rows(17).AutoFit()

# --- Row 18: label becomes "Método:", value becomes the Herlandí bio line ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Rows(18).RowHeight = 60

# --- Row 19: label becomes "Critério:", value becomes the "Aulas expositivas" text ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Range("C19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Rows(19).RowHeight = 60

# --- Row 20: label becomes "Norma de recuperação:", value becomes the "Média Aritmética" text ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"
$ws.Range("C20").Value = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"
$ws.Rows(20).RowHeight = 60

# --- Row 21: label becomes "Bibliografia:", value becomes the "NF = (MF+PR)/2" text ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação"
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação"
$ws.Rows(21).RowHeight = 120

# --- Row 22: the old bibliography row is removed entirely, shrinking the
#     sheet from 22 to 21 rows ---
$ws.Rows(22).Delete()
